$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header columns: "location" -> "name", "network" -> "netid"
$ws.Range("B2").Value = "name"
$ws.Range("C2").Value = "netid"

# Move active selection to B3 (matches author's final cursor position)
$ws.Range("B3").Select()
